# [Kadastro App] Yeni kayit eklendi: 3005
#
# Appends the new record (Kayit No 3005, Erdemli, 2025-09-11) as the next
# row after the existing data on both the master "Kayitlar" sheet and the
# "Erdemli" filtered sheet (sheet1 / sheet8), mirroring how the previous
# rows (e.g. row 64) are laid out.

$wb = $excel.ActiveWorkbook

$newRecord = @{
    KayitNo    = "3005"
    Tarih      = "2025-09-11"
    Birim      = "Erdemli"
    ParselSayisi = "1"
    Is         = "ÇAP"
    Personeller = "AYHAN KARADAYI (K.Teknisyeni)"
}

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the first empty row right after the existing data (row 65 here).
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    # Columns A, B and D hold numeric-looking text ("3005", "2025-09-11",
    # "1") just like the rest of the sheet's "number stored as text" data,
    # so force Text formatting before writing them to keep them as strings
    # instead of being auto-converted to a number/date.
    $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 2)).NumberFormat = "@"
    $ws.Cells.Item($newRow, 4).NumberFormat = "@"

    $ws.Cells.Item($newRow, 1).Value = $newRecord.KayitNo
    $ws.Cells.Item($newRow, 2).Value = $newRecord.Tarih
    $ws.Cells.Item($newRow, 3).Value = $newRecord.Birim
    $ws.Cells.Item($newRow, 4).Value = $newRecord.ParselSayisi
    $ws.Cells.Item($newRow, 5).Value = $newRecord.Is
    $ws.Cells.Item($newRow, 6).Value = $newRecord.Personeller

    # Keep Excel's "number stored as text" warning suppressed for the new
    # row, same as the rest of the column (mirrors the ignoredErrors/
    # numberStoredAsText range growing to include the new row).
    $newRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 6))
    try {
        $newRange.Errors.Item(9).Ignore = $true
    } catch {
    }
}

Write-Host "Added record 3005 to sheets: $($sheetNames -join ', ')"
